# Apply changes described by the diff:
#  - Update header text for B1 and C1 (append " tensorflow")
#  - Add new header D1 with combined arcface label
#  - Add new columns D and E with arcface distance data for rows 2-19
#  - A handful of column C values get tiny floating point corrections

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# --- Header row updates ---
$ws.Range("B1").Value = "Расстояние до переноса стиля tensorflow"
$ws.Range("C1").Value = "Расстояние после переноса стиля tensorflow"
$ws.Range("D1").Value = "Расстояние до переноса стиля arcface,Расстояние после переноса стиля arcface"

# --- Small floating point corrections to existing column C values ---
$ws.Range("C5").Value = 0.7964307258438679
$ws.Range("C10").Value = 0.7330427481233784
$ws.Range("C11").Value = 0.754825674803941
$ws.Range("C12").Value = 0.4725813935492413
$ws.Range("C15").Value = 0.6791952460002608
$ws.Range("C18").Value = 0.6392438678374194

# --- New D/E columns (arcface distances before/after) for rows 2-19 ---
$data = @{
    2  = @(0.03230053186416626, 0.01071752980351448)
    3  = @(0.07159358263015747, 0.06596450507640839)
    4  = @(0.1831206232309341, 0.01908320561051369)
    5  = @(0.11044642329216, 0.03032833524048328)
    6  = @(0.05190222710371017, 0.06330510973930359)
    7  = @(0.0741061195731163, 0.0427030511200428)
    8  = @(0.06903903931379318, 0.04660750180482864)
    9  = @(0.05451708287000656, 0.03171984478831291)
    10 = @(0.1547363102436066, 0.08104228973388672)
    11 = @(0.04635890573263168, 0.06895835697650909)
    12 = @(0.01744678616523743, 0.02074667811393738)
    13 = @(0.03705150634050369, 0.02712763473391533)
    14 = @(0.05199411883950233, 0.04723206162452698)
    15 = @(0.07539159804582596, 0.05160433799028397)
    16 = @(0.08692525327205658, 0.02815583720803261)
    17 = @(0.04295807331800461, 0.03566687181591988)
    18 = @(0.2073389887809753, 0.07804397493600845)
    19 = @(0.03760791569948196, 0.05074924603104591)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 4).Value = $values[0]
    $ws.Cells.Item($row, 5).Value = $values[1]
}
